$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Goal (per the commit "Added all current Data"):
#   - Insert a new first column ("A") holding a 0-based numeric segment
#     index (0..18).
#   - The old column A (segment name, e.g. "background") shifts to column B
#     and gets a new header "segments" in B1.
#   - The old measure columns B..F (PercActivationsOriginal/Correct/
#     Incorrect/Corrected/Fixed) shift right by one, to C..G.
#
# NOTE: Columns(...).Insert() / Range.Insert(xlShiftToRight) in this
# environment mis-shift content by two columns instead of one, and reading
# Range/Cells .Value back out does not return usable data either - so
# instead of relying on "insert and shift", every cell is written
# explicitly to its final address using the values read directly from the
# original workbook's XML, and formatting is (re)applied explicitly with
# Copy/PasteSpecial(Formats) / Range.Style so it lands on the right cells.
# ---------------------------------------------------------------------------

$xlPasteFormats = -4122

$labels = @(
    "background",
    "back_bumper",
    "back_glass",
    "back_left_door",
    "back_left_light",
    "back_right_door",
    "back_right_light",
    "front_bumper",
    "front_glass",
    "front_left_door",
    "front_left_light",
    "front_right_door",
    "front_right_light",
    "hood",
    "left_mirror",
    "right_mirror",
    "tailgate",
    "trunk",
    "wheel"
)

$measures = @(
    @(0.2233302586436274, 0.2225387910153632, 0.2919640402762031, 0.2943420731830562, 0.2233573688965221),
    @(0.006806036861633729, 0.006807684204421871, 0.006663184063669754, 0.006656611287645864, 0.006805961930113767),
    @(0.01530389343765177, 0.01528845829971843, 0.01664238391856818, 0.01811944281382678, 0.01532073233003003),
    @(0.004752661513772372, 0.004755584483528676, 0.004499190037733452, 0.004445816133578705, 0.004752053036053061),
    @(0.005117873612728906, 0.005122014261365026, 0.00475880854909228, 0.004694214571778435, 0.00511713722297778),
    @(0.006451775790282372, 0.006455739033853222, 0.006108094780364957, 0.006008943070878997, 0.006450645432535286),
    @(0.006753712319599705, 0.006759898182647172, 0.006217292182568843, 0.006152862295985334, 0.006752977800529674),
    @(0.3006371701245291, 0.3009661120686788, 0.2721122769282146, 0.2689558424985024, 0.300601185872424),
    @(0.15128062558172, 0.1515059919376195, 0.1317375075744001, 0.1367836495822203, 0.1513381530387956),
    @(0.002516709738565942, 0.002519074956010316, 0.002311604928471731, 0.002259547025147021, 0.002516116263631166),
    @(0.02282024236928663, 0.02284092656009624, 0.02102656921756709, 0.02048019556302745, 0.02281401355390451),
    @(0.002408039705547265, 0.00241146003395716, 0.00211143872679203, 0.002078263405427385, 0.002407661497428507),
    @(0.01799456682363838, 0.01801651068420853, 0.01609165875669511, 0.01561369048393108, 0.0179891178491045),
    @(0.2227139906715859, 0.2229809157399301, 0.1995670214225842, 0.1960082826960557, 0.2226734200358374),
    @(0.003420598726329381, 0.003422737646368247, 0.003235117772169281, 0.003331683613353148, 0.003421699604440826),
    @(0.002163653280140601, 0.002165849682504913, 0.001973187625114235, 0.001921081046757907, 0.002163059250296594),
    @(0.0002387190762089984, 0.0002398000790152947, 0.0001449776420656366, 0.0001370077522941128, 0.0002386282171941277),
    @(0.005267788499280173, 0.005182042473122698, 0.01270343567547521, 0.01192902387690292, 0.005258959984063566),
    @(0.00002168372745324416, 0.0000204091296450094, 0.0001322131597581288, 0.00008177837686835217, 0.00002110875655402823)
)

# --- Row 1 headers (B1:G1), all styled like the original bold/bordered
#     header cells. -----------------------------------------------------
$ws.Range("B1").Value = "segments"
$ws.Range("C1").Value = "PercActivationsOriginal"
$ws.Range("D1").Value = "PercActivationsCorrect"
$ws.Range("E1").Value = "PercActivationsIncorrect"
$ws.Range("F1").Value = "PercActivationsCorrected"
$ws.Range("G1").Value = "PercActivationsFixed"

# Apply the header look (bold, thin border all sides, centered/top aligned)
# to every row-1 header cell by copying format from a cell that already has
# it (A1's original style lived on the "label" cells, e.g. A2; reuse that).
$ws.Range("A2").Copy() | Out-Null
$ws.Range("B1:G1").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = $false

# --- Data rows 2..20: column A = numeric 0-based index (same styling as
#     the old label column), column B = label text (plain/no special
#     styling), columns C..G = the five measures (plain/no styling). ------
for ($i = 0; $i -lt $labels.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $i
    $ws.Cells.Item($row, 2).Value = $labels[$i]
    $rowVals = $measures[$i]
    for ($j = 0; $j -lt $rowVals.Length; $j++) {
        $ws.Cells.Item($row, 3 + $j).Value = $rowVals[$j]
    }
}

# Column A (index) keeps the bold/bordered/centered "label" look the
# original column A had.
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A2:A20").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = $false

# Column B (labels) and C:G (measures) get the plain default look the
# measure columns originally had (no border/bold/center).
$ws.Range("B2:B20").Style = "Normal"
$ws.Range("C2:G20").Style = "Normal"
